# subject_template.xlsx: add a "studentCount" column between subName/sem,
# fill in its value (60) for the sample row, and tidy up the credits cell
# / hyperlink that end up shifting one column to the right as a result.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Shift columns E..H (both header + data row) one column to the right,
# opening up column E for the new "studentCount" field. We shift cell by
# cell (right to left, so we never clobber a cell before it has moved) via
# Cut/Paste rather than Columns.Insert() because a plain Cut carries the
# original cell style along with it, whereas Insert() leaves the worksheet's
# <hyperlink> anchor pointing at the old (pre-shift) address.
$ws.Range("H1").Cut($ws.Range("I1"))
$ws.Range("G1").Cut($ws.Range("H1"))
$ws.Range("F1").Cut($ws.Range("G1"))
$ws.Range("E1").Cut($ws.Range("F1"))
$ws.Range("H2").Cut($ws.Range("I2"))
$ws.Range("G2").Cut($ws.Range("H2"))
$ws.Range("F2").Cut($ws.Range("G2"))
$ws.Range("E2").Cut($ws.Range("F2"))

# The now-empty column E should look like the rest of the plain data
# columns (same formatting as column D) before we put the new header/value
# into it.
$ws.Range("D1").Copy()
$ws.Range("E1").PasteSpecial(-4122)
$ws.Range("D2").Copy()
$ws.Range("E2").PasteSpecial(-4122)

$ws.Range("E1").Value = "studentCount"
$ws.Range("E2").Value = 60

# credits value (kept as a plain integer, same as the source workbook)
$ws.Range("I2").Value = 3

# --- Fix up the hyperlink, which lived on the old "degree" cell (F2) and
# has now moved to G2. Stash the cell's current (correct) formatting first,
# recreate the hyperlink (Excel always reformats the anchor cell with its
# built-in "Hyperlink" style when a link is (re)created) and then restore
# the original look so the cell keeps matching the rest of the sheet.
$ws.Range("G2").Copy()
$ws.Range("Z1").PasteSpecial(-4122)
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("G2"), "http://m.sc/")
$ws.Range("Z1").Copy()
$ws.Range("G2").PasteSpecial(-4122)
$ws.Range("Z1").Clear()

# Normalize the rich-text "M.Sc. " label back to 10pt (it previously had no
# explicit size, which round-trips inconsistently once the workbook is
# resaved).
$ws.Range("G2").Characters(1, 4).Font.Size = 10
$ws.Range("G2").Characters(5, 2).Font.Size = 10

# Match the author's on-disk cursor position.
[void]$ws.Range("G15").Select()
